# Updated BoM with new part for C1Mx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 corresponds to the C1Mx line item (designator in I7).
# Replace the old Vishay/VJ0603... ceramic cap part with the new Murata/GRM part.
$ws.Range("B7").Value = "SMD Multilayer Ceramic Capacitor, GRM Series, 0.01 µF, ± 10%, X7R, 50 V, 0603 [1608 Metric]"
$ws.Range("D7").Value = "MURATA"
$ws.Range("E7").Value = "GRM188R71H103KA01D"
$ws.Range("G7").Value = "38K1669"
$ws.Range("H7").Value = 0.008

# Move the active selection to H8, matching the saved view state.
$ws.Range("H8").Select()
